$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 852 (shifts the existing row 852..893 down to 853..894)
$ws.Rows.Item(852).Insert()

# Column A holds the date as literal text (e.g. "2026/12/29"), not a real date
# value, so force the cell to text formatting before writing it to stop Excel
# from auto-converting the "yyyy/mm/dd"-looking string into a date serial.
$ws.Range("A852").NumberFormat = "@"
$ws.Range("A852").Value = "2026/02/25"
# Reset back to the default "Normal" style so the cell matches its neighbours
# (no explicit style index), now that the text value is already locked in.
$ws.Range("A852").Style = "Normal"

$ws.Range("B852").Value = "水"
$ws.Range("C852").Value = 1
$ws.Range("D852").Value = 48
